$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2 through 453) holds a "last changed" date serial number.
# Update it from 45204 (2023-10-05) to 45205 (2023-10-06) for every data row.
$ws.Range("C2:C453").Value = 45205
